$d = $word.ActiveDocument

$d.Content.Find.Execute("321÷8=40, 1", $true, $false, $false, $false, $false, $true, 1, $false, "228÷4=57, 0", 2) | Out-Null
$d.Content.Find.Execute("494÷4=123, 2", $true, $false, $false, $false, $false, $true, 1, $false, "165÷4=41, 1", 2) | Out-Null
$d.Content.Find.Execute("357÷6=59, 3", $true, $false, $false, $false, $false, $true, 1, $false, "611÷3=203, 2", 2) | Out-Null
$d.Content.Find.Execute("380÷8=47, 4", $true, $false, $false, $false, $false, $true, 1, $false, "245÷4=61, 1", 2) | Out-Null
$d.Content.Find.Execute("788÷3=262, 2", $true, $false, $false, $false, $false, $true, 1, $false, "154÷5=30, 4", 2) | Out-Null
$d.Content.Find.Execute("737÷4=184, 1", $true, $false, $false, $false, $false, $true, 1, $false, "531÷6=88, 3", 2) | Out-Null
$d.Content.Find.Execute("155÷2=77, 1", $true, $false, $false, $false, $false, $true, 1, $false, "301÷8=37, 5", 2) | Out-Null
$d.Content.Find.Execute("425÷5=85, 0", $true, $false, $false, $false, $false, $true, 1, $false, "891÷2=445, 1", 2) | Out-Null
$d.Content.Find.Execute("586÷2=293, 0", $true, $false, $false, $false, $false, $true, 1, $false, "163÷7=23, 2", 2) | Out-Null
$d.Content.Find.Execute("139÷2=69, 1", $true, $false, $false, $false, $false, $true, 1, $false, "893÷6=148, 5", 2) | Out-Null
$d.Content.Find.Execute("186÷4=46, 2", $true, $false, $false, $false, $false, $true, 1, $false, "436÷2=218, 0", 2) | Out-Null
$d.Content.Find.Execute("565÷7=80, 5", $true, $false, $false, $false, $false, $true, 1, $false, "498÷8=62, 2", 2) | Out-Null
$d.Content.Find.Execute("979÷3=326, 1", $true, $false, $false, $false, $false, $true, 1, $false, "415÷7=59, 2", 2) | Out-Null
$d.Content.Find.Execute("949÷4=237, 1", $true, $false, $false, $false, $false, $true, 1, $false, "216÷4=54, 0", 2) | Out-Null
$d.Content.Find.Execute("607÷4=151, 3", $true, $false, $false, $false, $false, $true, 1, $false, "119÷5=23, 4", 2) | Out-Null
$d.Content.Find.Execute("491÷5=98, 1", $true, $false, $false, $false, $false, $true, 1, $false, "794÷4=198, 2", 2) | Out-Null
$d.Content.Find.Execute("680÷3=226, 2", $true, $false, $false, $false, $false, $true, 1, $false, "627÷6=104, 3", 2) | Out-Null
$d.Content.Find.Execute("557÷8=69, 5", $true, $false, $false, $false, $false, $true, 1, $false, "327÷4=81, 3", 2) | Out-Null
$d.Content.Find.Execute("409÷7=58, 3", $true, $false, $false, $false, $false, $true, 1, $false, "337÷7=48, 1", 2) | Out-Null
$d.Content.Find.Execute("396÷4=99, 0", $true, $false, $false, $false, $false, $true, 1, $false, "780÷7=111, 3", 2) | Out-Null
$d.Content.Find.Execute("729÷7=104, 1", $true, $false, $false, $false, $false, $true, 1, $false, "716÷2=358, 0", 2) | Out-Null
$d.Content.Find.Execute("620÷8=77, 4", $true, $false, $false, $false, $false, $true, 1, $false, "630÷7=90, 0", 2) | Out-Null
$d.Content.Find.Execute("501÷3=167, 0", $true, $false, $false, $false, $false, $true, 1, $false, "562÷9=62, 4", 2) | Out-Null
$d.Content.Find.Execute("142÷7=20, 2", $true, $false, $false, $false, $false, $true, 1, $false, "301÷7=43, 0", 2) | Out-Null
$d.Content.Find.Execute("624÷6=104, 0", $true, $false, $false, $false, $false, $true, 1, $false, "425÷4=106, 1", 2) | Out-Null
